# Update "paises" (countries) covid tracking sheet:
#  - refresh the "Datos actualizados" timestamp
#  - refresh case counters for a set of countries
#  - Marruecos overtakes Paises Bajos / Emiratos Arabes Unidos in the ranking
#  - Etiopia overtakes Portugal in the ranking

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 20:53"

# --- Row 4: Estados Unidos ----------------------------------------------
$ws.Range("B4").Value = 6529319
$ws.Range("C4").Value = 15088
$ws.Range("D4").Value = 3817379
$ws.Range("E4").Value = 2517316
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 594
$ws.Range("H4").Value = 194624

# --- Row 5: India ---------------------------------------------------------
$ws.Range("B5").Value = 4462965
$ws.Range("C5").Value = 95529
$ws.Range("D5").Value = 3466819
$ws.Range("E5").Value = 921055
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1168
$ws.Range("H5").Value = 75091

# --- Row 17: Francia -------------------------------------------------------
$ws.Range("B17").Value = 344101
$ws.Range("C17").Value = 8577
$ws.Range("D17").Value = 88226
$ws.Range("E17").Value = 225081
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 30794

# --- Row 29: Canada ----------------------------------------------------
$ws.Range("B29").Value = 134093
$ws.Range("C29").Value = 345
$ws.Range("D29").Value = 117945
$ws.Range("E29").Value = 6994
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 9154

# --- Rows 44-46: Marruecos jumps ahead of Paises Bajos & Emiratos Arabes Unidos
# Row 44 becomes Marruecos (new, updated data)
$ws.Range("A44").Value = "Marruecos"
$ws.Range("B44").Value = 77878
$ws.Range("C44").Value = 2157
$ws.Range("D44").Value = 59723
$ws.Range("E44").Value = 16702
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 26
$ws.Range("H44").Value = 1453

# Row 45 becomes Paises Bajos (its data is unchanged, just shifted down a rank)
$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("B45").Value = 77688
$ws.Range("C45").Value = 1140
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 6246

# Row 46 becomes Emiratos Arabes Unidos (its data is unchanged, just shifted down a rank)
$ws.Range("A46").Value = "Emiratos Arabes Unidos"
$ws.Range("B46").Value = 75981
$ws.Range("C46").Value = 883
$ws.Range("D46").Value = 67359
$ws.Range("E46").Value = 8229
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 393

# --- Rows 51-52: Etiopia jumps ahead of Portugal ---------------------------
# Row 51 becomes Etiopia (new, updated data)
$ws.Range("A51").Value = "Etiopia"
$ws.Range("B51").Value = 61700
$ws.Range("C51").Value = 916
$ws.Range("D51").Value = 23054
$ws.Range("E51").Value = 37680
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 17
$ws.Range("H51").Value = 966

# Row 52 becomes Portugal (its data is unchanged, just shifted down a rank)
$ws.Range("A52").Value = "Portugal"
$ws.Range("B52").Value = 61541
$ws.Range("C52").Value = 646
$ws.Range("D52").Value = 43284
$ws.Range("E52").Value = 16408
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 1849

# --- Row 54: Barein ----------------------------------------------------
$ws.Range("B54").Value = 56778
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 51574
$ws.Range("E54").Value = 5001
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 203

# --- Row 101: Maldivas ---------------------------------------------------
$ws.Range("B101").Value = 8834
$ws.Range("C101").Value = 93
$ws.Range("D101").Value = 6288
$ws.Range("E101").Value = 2515
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 31

# --- Row 149: Yemen --------------------------------------------------------
$ws.Range("B149").Value = 1999
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 1209
$ws.Range("E149").Value = 214
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 576
